$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: bump issue number and report week dates ---
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"

# --- Column width fixes: columns E and H should match the narrower
#     standard width used by the other numeric columns (F, G, etc.) ---
$ws.Columns.Item(5).ColumnWidth = 5.43
$ws.Columns.Item(8).ColumnWidth = 5.43

# --- Weekly crime-statistics figures (rows 14-28) ---
    # Row 14
    $ws.Range("F14").Value = 1
    $ws.Range("M14").Value = 100
    $ws.Range("M14").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("N14").Value = -50

    # Row 15
    $ws.Range("C15").Value = 2
    $ws.Range("C15").NumberFormat = '#,##0'
    $ws.Range("D15").Value = 2
    $ws.Range("E15").Value = 0
    $ws.Range("F15").Value = 4
    $ws.Range("G15").Value = 4
    $ws.Range("H15").Value = 0
    $ws.Range("I15").Value = 9
    $ws.Range("J15").Value = 6
    $ws.Range("K15").Value = 50
    $ws.Range("L15").Value = 125
    $ws.Range("M15").Value = 125
    $ws.Range("N15").Value = 50

    # Row 16
    $ws.Range("G16").Value = 9
    $ws.Range("H16").Value = -33.333333333333
    $ws.Range("I16").Value = 15
    $ws.Range("J16").Value = 30
    $ws.Range("K16").Value = -50
    $ws.Range("L16").Value = -46.428571428571
    $ws.Range("M16").Value = -55.882352941176
    $ws.Range("N16").Value = -87.903225806451

    # Row 17
    $ws.Range("C17").Value = 9
    $ws.Range("E17").Value = 28.571428571428
    $ws.Range("F17").Value = 34
    $ws.Range("H17").Value = 9.677419354838
    $ws.Range("I17").Value = 64
    $ws.Range("J17").Value = 73
    $ws.Range("K17").Value = -12.328767123287
    $ws.Range("L17").Value = 56.097560975609
    $ws.Range("M17").Value = 68.421052631578
    $ws.Range("N17").Value = -28.089887640449

    # Row 18
    $ws.Range("C18").Value = 2
    $ws.Range("C18").NumberFormat = '#,##0'
    $ws.Range("D18").Value = 1
    $ws.Range("E18").Value = 100
    $ws.Range("F18").Value = 6
    $ws.Range("G18").Value = 7
    $ws.Range("H18").Value = -14.285714285714
    $ws.Range("I18").Value = 11
    $ws.Range("J18").Value = 18
    $ws.Range("K18").Value = -38.888888888888
    $ws.Range("L18").Value = -38.888888888888
    $ws.Range("M18").Value = -56
    $ws.Range("N18").Value = -90.434782608695

    # Row 19
    $ws.Range("C19").Value = 2
    $ws.Range("D19").Value = 1
    $ws.Range("D19").NumberFormat = '#,##0'
    $ws.Range("E19").Value = 100
    $ws.Range("E19").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("G19").Value = 8
    $ws.Range("H19").Value = 50
    $ws.Range("I19").Value = 30
    $ws.Range("J19").Value = 33
    $ws.Range("K19").Value = -9.090909090909
    $ws.Range("L19").Value = -23.076923076923
    $ws.Range("M19").Value = 42.857142857142
    $ws.Range("N19").Value = -49.152542372881

    # Row 20
    $ws.Range("C20").Value = 2
    $ws.Range("D20").Value = 1
    $ws.Range("E20").Value = 100
    $ws.Range("F20").Value = 3
    $ws.Range("G20").Value = 7
    $ws.Range("H20").Value = -57.142857142857
    $ws.Range("I20").Value = 6
    $ws.Range("J20").Value = 14
    $ws.Range("K20").Value = -57.142857142857
    $ws.Range("L20").Value = -62.5
    $ws.Range("M20").Value = -71.428571428571
    $ws.Range("N20").Value = -94.594594594594

    # Row 21
    $ws.Range("C21").Value = 19
    $ws.Range("D21").Value = 14
    $ws.Range("E21").Value = 35.714285714285
    $ws.Range("F21").Value = 66
    $ws.Range("G21").Value = 66
    $ws.Range("H21").Value = 0
    $ws.Range("I21").Value = 137
    $ws.Range("J21").Value = 175
    $ws.Range("K21").Value = -21.714285714285
    $ws.Range("L21").Value = -6.164383561643
    $ws.Range("M21").Value = -4.861111111111
    $ws.Range("N21").Value = -73.031496062992

    # Row 23
    $ws.Range("C23").Value = 3
    $ws.Range("E23").Value = 50
    $ws.Range("F23").Value = 8
    $ws.Range("G23").Value = 6
    $ws.Range("H23").Value = 33.333333333333
    $ws.Range("I23").Value = 15
    $ws.Range("J23").Value = 19
    $ws.Range("K23").Value = -21.052631578947
    $ws.Range("L23").Value = -21.052631578947
    $ws.Range("M23").Value = 25

    # Row 24
    $ws.Range("C24").Value = 16
    $ws.Range("D24").Value = 6
    $ws.Range("E24").Value = 166.666666666667
    $ws.Range("F24").Value = 53
    $ws.Range("G24").Value = 41
    $ws.Range("H24").Value = 29.268292682926
    $ws.Range("I24").Value = 154
    $ws.Range("J24").Value = 135
    $ws.Range("K24").Value = 14.074074074074
    $ws.Range("L24").Value = 35.087719298245
    $ws.Range("M24").Value = 123.188405797101

    # Row 25
    $ws.Range("C25").Value = 5
    $ws.Range("C25").NumberFormat = '#,##0'
    $ws.Range("D25").Value = 1
    $ws.Range("D25").NumberFormat = '#,##0'
    $ws.Range("E25").Value = 400
    $ws.Range("E25").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("F25").Value = 11
    $ws.Range("G25").Value = 4
    $ws.Range("H25").Value = 175
    $ws.Range("I25").Value = 22
    $ws.Range("J25").Value = 10
    $ws.Range("K25").Value = 120
    $ws.Range("L25").Value = 46.666666666666

    # Row 26
    $ws.Range("D26").Value = 4
    $ws.Range("E26").Value = 125
    $ws.Range("F26").Value = 37
    $ws.Range("G26").Value = 22
    $ws.Range("H26").Value = 68.181818181818
    $ws.Range("I26").Value = 104
    $ws.Range("J26").Value = 80
    $ws.Range("K26").Value = 30
    $ws.Range("L26").Value = 30
    $ws.Range("M26").Value = 28.395061728395

    # Row 27
    $ws.Range("C27").Value = 2
    $ws.Range("C27").NumberFormat = '#,##0'
    $ws.Range("E27").Value = 0
    $ws.Range("F27").Value = 4
    $ws.Range("G27").Value = 5
    $ws.Range("H27").Value = -20
    $ws.Range("I27").Value = 10
    $ws.Range("J27").Value = 7
    $ws.Range("K27").Value = 42.857142857142
    $ws.Range("L27").Value = 0

    # Row 28
    $ws.Range("F28").Value = 1
    $ws.Range("G28").Value = 1
    $ws.Range("L28").Value = -40


